# Competition.xlsx update: SNOW, JNJ, PFE updates
# - Main sheet: append the new day's NAV row (row 12)
# - Trades sheet: append the new day's trade log (rows 217-223)

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) "Main" sheet - add row 12 values/formulas (new NAV datapoint)
# ----------------------------------------------------------------------
$main = $wb.Worksheets.Item("Main")

$main.Range("C12").Value = 53440.36
$main.Range("C12").NumberFormat = "#,##0.00"
$main.Range("C12").HorizontalAlignment = -4152

$main.Range("D12").Formula = '=C12-C11'
$main.Range("D12").NumberFormat = "#,##0.00"
$main.Range("D12").HorizontalAlignment = -4152

$main.Range("E12").Formula = '=C12-$C$4'
$main.Range("E12").NumberFormat = "#,##0"
$main.Range("E12").HorizontalAlignment = -4152

$main.Range("F12").Formula = '=+C12/C11-1'
$main.Range("F12").NumberFormat = "0.0%"
$main.Range("F12").HorizontalAlignment = -4152

$main.Range("G12").Formula = '=C12/$C$4-1'
$main.Range("G12").NumberFormat = "0.0%"

$main.Range("C12").Select()

# ----------------------------------------------------------------------
# 2) "Trades" sheet - append new day's trade rows (217-223)
# ----------------------------------------------------------------------
$trades = $wb.Worksheets.Item("Trades")

# Row 217: MSFT trade
$trades.Range("B217").Value = "MSFT"
$trades.Range("C217").Value = "2022-07-27, 09:37:29"
$trades.Range("D217").Value = -20
$trades.Range("E217").Value = 264.17
$trades.Range("F217").Value = 268.74
$trades.Range("G217").Value = 5283.4
$trades.Range("H217").Value = -1.1200000000000001
$trades.Range("I217").Value = -5282.28
$trades.Range("J217").Value = 10.73
$trades.Range("K217").Value = -91.4
$trades.Range("L217").Value = "O"

# Row 218: Total MSFT
$trades.Range("B218").Value = "Total MSFT"
$trades.Range("D218").Value = -20
$trades.Range("E218").Value = " "
$trades.Range("G218").Value = 5283.4
$trades.Range("H218").Value = -1.1200000000000001
$trades.Range("I218").Value = -5282.28
$trades.Range("J218").Value = 10.73
$trades.Range("K218").Value = -91.4
$trades.Range("L218").Value = " "

# Row 219: SGHC trade
$trades.Range("B219").Value = "SGHC"
$trades.Range("C219").Value = "2022-07-27, 14:41:05"
$trades.Range("D219").Value = -500
$trades.Range("E219").Value = 4.0999999999999996
$trades.Range("F219").Value = 3.94
$trades.Range("G219").Value = 2050
$trades.Range("H219").Value = -2.61
$trades.Range("I219").Value = -2047.39
$trades.Range("J219").Value = 0
$trades.Range("K219").Value = 80
$trades.Range("L219").Value = "O"

# Row 220: Total SGHC
$trades.Range("B220").Value = "Total SGHC"
$trades.Range("D220").Value = -500
$trades.Range("E220").Value = " "
$trades.Range("G220").Value = 2050
$trades.Range("H220").Value = -2.61
$trades.Range("I220").Value = -2047.39
$trades.Range("J220").Value = 0
$trades.Range("K220").Value = 80
$trades.Range("L220").Value = " "

# Row 221: TEAM trade
$trades.Range("B221").Value = "TEAM"
$trades.Range("C221").Value = "2022-07-27, 10:44:44"
$trades.Range("D221").Value = 10
$trades.Range("E221").Value = 191
$trades.Range("F221").Value = 197.28
$trades.Range("G221").Value = -1910
$trades.Range("H221").Value = -1
$trades.Range("I221").Value = 1911
$trades.Range("J221").Value = 0
$trades.Range("K221").Value = 62.8
$trades.Range("L221").Value = "O"

# Row 222: Total TEAM
$trades.Range("B222").Value = "Total TEAM"
$trades.Range("D222").Value = 10
$trades.Range("E222").Value = " "
$trades.Range("G222").Value = -1910
$trades.Range("H222").Value = -1
$trades.Range("I222").Value = 1911
$trades.Range("J222").Value = 0
$trades.Range("K222").Value = 62.8
$trades.Range("L222").Value = " "

# Row 223: grand Total (day total, recalculated)
$trades.Range("B223").Value = "Total"
$trades.Range("G223").Value = 5423.4
$trades.Range("H223").Value = -4.74
$trades.Range("I223").Value = -5418.66
$trades.Range("J223").Value = 10.73
$trades.Range("K223").Value = 51.4
$trades.Range("L223").Value = " "

# Apply number formatting (style 6 -> numFmt "#,##0.00", no alignment) to the
# numeric cells of the new rows that use that style in the template rows above.
$fmtCells = @(
  "E217","F217","G217","H217","I217","J217","K217",
  "E218","G218","H218","I218","J218","K218",
  "E219","F219","G219","H219","I219","J219","K219",
  "E220","G220","H220","I220","J220","K220",
  "E221","F221","G221","H221","I221","J221","K221",
  "E222","G222","H222","I222","J222","K222",
  "G223","I223","J223","K223"
)
# NB: H223 deliberately excluded - it keeps the default (unstyled) format,
# matching the existing "Total" row pattern (see H215 in the template).
foreach ($addr in $fmtCells) {
  $trades.Range($addr).NumberFormat = "#,##0.00"
}

# L218/L220/L222/L223/E218/E220/E222 are text placeholders (" ") drawn with
# the same numeric style (s="6") as their row's other numeric cells.
foreach ($addr in @("L218","L220","L222","L223")) {
  $trades.Range($addr).NumberFormat = "#,##0.00"
}

$trades.Range("G223:L223").Select()
